# Add the "Standard Kalk Kat3" product as a new column (K) to the lime
# products table, mirroring the "Miljøkalk VK3" (column D) values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header
$ws.Range("K1").Value = "Standard Kalk Kat3"

# Data rows (2-15) - same values as column D (Miljøkalk VK3)
$ws.Range("K2").Value = 39.6
$ws.Range("K3").Value = 0.4
$ws.Range("K4").Value = 0.7
$ws.Range("K5").Value = 2
$ws.Range("K6").Value = 61.2
$ws.Range("K7").Value = 60.5
$ws.Range("K8").Value = 60.2
$ws.Range("K9").Value = 51.7
$ws.Range("K10").Value = 48.5
$ws.Range("K11").Value = 1
$ws.Range("K12").Value = 1.2
$ws.Range("K13").Value = 1.6
$ws.Range("K14").Value = 2.2
$ws.Range("K15").Value = 3.2

# Resize column K to fit its (wider) header text
$ws.Columns.Item(11).ColumnWidth = 16.6

# Leave the selection where the author's session ended up
$ws.Range("K19").Select() | Out-Null
